$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.377412438392639
$ws.Range("B1").Value = 2.652031421661377
$ws.Range("C1").Value = 5.805872440338135
$ws.Range("D1").Value = 2.257520437240601
$ws.Range("E1").Value = 1.214731216430664
